# chore(runtime): publish files + archive (2025-12-01 11:05:52)
# Appends the 2025-11-30 KHL match results to Matches_SOG, refreshes the
# as_of_utc timestamp + derived shot-on-goal aggregates on Shots_HA /
# Shots_Summary, and bumps Meta_ext's as_of_utc + build_version.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Matches_SOG: append three new match rows (338-340)
# ---------------------------------------------------------------------
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

$newMatches = @(
    @("897831", "2025-11-30T10:00:00", "Адмирал", "Динамо Мн", 31, 38, "khl_text"),
    @("897832", "2025-11-30T10:00:00", "Амур",    "ХК Сочи",   39, 21, "khl_text"),
    @("897833", "2025-11-30T17:00:00", "Ак Барс", "Драконы",   32, 36, "khl_text")
)

$row = 338
foreach ($m in $newMatches) {
    # uid looks numeric ("897831") but is stored as text in the source
    # sheet, matching every other row. Force text storage via a "@"
    # number format (then restore the default "Normal" style so no
    # stray style index is left on the cell) instead of a leading
    # apostrophe, which would also mark the cell quotePrefix="1".
    $wsMatches.Cells.Item($row, 1).NumberFormat = "@"
    $wsMatches.Cells.Item($row, 1).Value = $m[0]
    $wsMatches.Cells.Item($row, 1).Style = "Normal"
    $wsMatches.Cells.Item($row, 2).Value = $m[1]
    $wsMatches.Cells.Item($row, 3).Value = $m[2]
    $wsMatches.Cells.Item($row, 4).Value = $m[3]
    $wsMatches.Cells.Item($row, 5).Value = $m[4]
    $wsMatches.Cells.Item($row, 6).Value = $m[5]
    $wsMatches.Cells.Item($row, 7).Value = $m[6]
    $row++
}

# ---------------------------------------------------------------------
# 2) Shots_HA: as_of_utc refresh (D2:D23) + updated HOG aggregates
# ---------------------------------------------------------------------
$wsHA = $wb.Worksheets.Item("Shots_HA")

$wsHA.Range("D2:D23").Value = "2025-11-30T17:00:00Z"

# row -> @{ col = newValue }
$haUpdates = @{
    4  = @{ E = 13; G = 487; H = 352; I = 37.5; J = 27.1 }
    5  = @{ E = 18; G = 599; H = 473; I = 33.3; J = 26.3 }
    6  = @{ E = 15; G = 458; H = 511; I = 30.5; J = 34.1 }
    9  = @{ F = 13; K = 450; L = 366; M = 34.6; N = 28.2 }
    10 = @{ F = 19; K = 531; L = 689; M = 27.9; N = 36.3 }
    22 = @{ F = 14; K = 369; L = 522; M = 26.4; N = 37.3 }
}

foreach ($r in $haUpdates.Keys) {
    foreach ($col in $haUpdates[$r].Keys) {
        $wsHA.Range("$col$r").Value = $haUpdates[$r][$col]
    }
}

# ---------------------------------------------------------------------
# 3) Shots_Summary: as_of_utc refresh (D2:D23) + updated SOG aggregates
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Shots_Summary")

$wsSummary.Range("D2:D23").Value = "2025-11-30T17:00:00Z"

$summaryUpdates = @{
    4  = @{ E = 29; F = 984;  G = 796;  H = 33.9; I = 27.4 }
    5  = @{ E = 33; F = 1109; G = 906;  H = 33.6; I = 27.5 }
    6  = @{ E = 31; F = 896;  G = 1125; H = 28.9; I = 36.3 }
    9  = @{ E = 31; F = 1111; G = 849;            I = 27.4 }
    10 = @{ E = 31; F = 873;  G = 1105; H = 28.2; I = 35.6 }
    22 = @{ E = 30; F = 841;  G = 1034; H = 28;   I = 34.5 }
}

foreach ($r in $summaryUpdates.Keys) {
    foreach ($col in $summaryUpdates[$r].Keys) {
        $wsSummary.Range("$col$r").Value = $summaryUpdates[$r][$col]
    }
}

# ---------------------------------------------------------------------
# 4) Meta_ext: as_of_utc + build_version bump
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Range("B2").Value = "2025-11-30T17:00:00Z"
$wsMeta.Range("D2").Value = 23
